# Commit: "Files must now only contain account, debit and credit column"
# This removes the extra "category" labels that used to live in column H
# (Bank Balances, Plant and Equipment, Trade Receivables, Trade Payables,
# GST Payables, Accruals, Amount owing ..., Current Income Tax Liabilities,
# Borrowings, Share Capital, Retained Profits, Revenue, Cost of Sales,
# Administrative / Distribution & Marketing / Finance Expenses, Other Income,
# Income Tax Expense, Deposits, ...) as well as the report title block text
# that used to sit in column E (VSIG Pte. Ltd. / Trial Balance / December 2016)
# so that the sheet only keeps the Account / Debit / Credit columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the title block (column E, rows 2-8) ---
# These keep their existing cell formatting, only the text is removed.
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("E8").Value = ""

# --- Clear the column-H "category" helper cells ---
# Rows 11-14 and the check-total row 69 carry a distinct style (border etc.)
# from the rest of the column, so only their content is removed while the
# style is preserved.
$ws.Range("H11").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("H13").Value = ""
$ws.Range("H14").Value = ""
$ws.Range("H69").Value = ""

# The remaining data rows in column H use the plain default column style,
# so fully clearing (content + formatting) drops the redundant cell.
$plainRows = @(16,17,18,19,20,21,22,23,24,25,27,28,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67)
foreach ($r in $plainRows) {
    $ws.Range("H$r").Clear()
}

# Rows 13 and 14 had an enlarged row height (24) to fit the two-line
# category text that lived in H13/H14; now that the text is gone they go
# back to the sheet's normal row height.
$ws.Rows("13").AutoFit()
$ws.Rows("14").AutoFit()

# Update the active selection to reflect where the edit ended up.
$ws.Range("H13").Select()
